$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number for every data row (2-130).
# The diff shows this value moving from 45204 (2023-10-05) to 45205 (2023-10-06)
# for every single row, so update the whole block in one shot.
$ws.Range("C2:C130").Value = 45205
